# Removing unneeded Word snippets
#
# The "Snippets" table contains two rows that reference snippet ids/method
# names that no longer exist ("addParagraphs" / "insertPageBreak"). This
# script removes those two whole rows from the worksheet, which shifts the
# remaining rows up and shrinks the backing table/range from A1:E35 to
# A1:E33.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Snippets")

# Row 17 (Paragraph / insertBreak / insertPageBreak) - delete first so the
# still-lower row 4 keeps its original row number until we remove it too.
$ws.Rows.Item(17).Delete()

# Row 4 (Body / insertParagraph / addParagraphs)
$ws.Rows.Item(4).Delete()

# Leave the selection where Excel would naturally land after this edit.
$ws.Range("Q10").Select()
